# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 20, pushing the existing
# rows 20-23 down to 21-24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 20 (shifts rows 20..23 -> 21..24)
$ws.Rows.Item(20).Insert()

# Populate the new row 20 with the latest weekly data
$ws.Cells.Item(20, 1).Value = 11
$ws.Cells.Item(20, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(20, 3).Value = "Bíobío"
$ws.Cells.Item(20, 4).Value = 45280
$ws.Cells.Item(20, 5).Value = 8
$ws.Cells.Item(20, 6).Value = "Fruta"
$ws.Cells.Item(20, 7).Value = 100101
$ws.Cells.Item(20, 8).Value = "Berries"
$ws.Cells.Item(20, 9).Value = 100101004
$ws.Cells.Item(20, 10).Value = "Frambuesa"
$ws.Cells.Item(20, 11).Value = "Sin especificar"
$ws.Cells.Item(20, 12).Value = "Primera"
$ws.Cells.Item(20, 13).Value = 80
$ws.Cells.Item(20, 14).Value = 6000
$ws.Cells.Item(20, 15).Value = 6000
$ws.Cells.Item(20, 16).Value = 6000
$ws.Cells.Item(20, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(20, 18).Value = "Provincia de Linares"
$ws.Cells.Item(20, 19).Value = 3000
$ws.Cells.Item(20, 20).Value = 2
